$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate data among rows 118-120 (odds data re-sequenced) ---
# Row 118
$ws.Range("A118").Value = 116
$ws.Range("B118").Value = 7013885
$ws.Range("C118").Value = "Uruguay Primera División"
$ws.Range("D118").Value = "Uruguay Clausura"
$ws.Range("E118").Value = 45267.70833333334
$ws.Range("F118").Value = "La Luz"
$ws.Range("G118").Value = "Atletico Fenix Montevideo"
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = 51
$ws.Range("K118").Value = 3
$ws.Range("L118").Value = 3
$ws.Range("M118").Value = 2.4
$ws.Range("N118").Value = 2.9
$ws.Range("O118").Value = 2.75
$ws.Range("P118").Value = 2.6
$ws.Range("Q118").Value = 0
$ws.Range("R118").Value = 2.025
$ws.Range("S118").Value = 1.825
$ws.Range("T118").Value = 2
$ws.Range("U118").Value = 2.025
$ws.Range("V118").Value = 1.825
$ws.Range("W118").Value = -1
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = 1.6
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = 0.825
$ws.Range("AB118").Value = 0
$ws.Range("AC118").Value = 0

# Row 119
$ws.Range("A119").Value = 117
$ws.Range("B119").Value = 7013409
$ws.Range("C119").Value = "Uruguay Primera División"
$ws.Range("D119").Value = "Uruguay Clausura"
$ws.Range("E119").Value = 45267.70833333334
$ws.Range("F119").Value = "Nacional De Football"
$ws.Range("G119").Value = "Torque"
$ws.Range("H119").Value = 1
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = 50
$ws.Range("K119").Value = 1.666
$ws.Range("L119").Value = 3.9
$ws.Range("M119").Value = 4.5
$ws.Range("N119").Value = 1.615
$ws.Range("O119").Value = 4
$ws.Range("P119").Value = 4.75
$ws.Range("Q119").Value = -0.75
$ws.Range("R119").Value = 1.8
$ws.Range("S119").Value = 2.05
$ws.Range("T119").Value = 2.75
$ws.Range("U119").Value = 1.95
$ws.Range("V119").Value = 1.9
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = 3
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = 1.05
$ws.Range("AB119").Value = -1
$ws.Range("AC119").Value = 0.8999999999999999

# Row 120
$ws.Range("A120").Value = 118
$ws.Range("B120").Value = 7013702
$ws.Range("C120").Value = "Uruguay Primera División"
$ws.Range("D120").Value = "Uruguay Clausura"
$ws.Range("E120").Value = 45267.70833333334
$ws.Range("F120").Value = "Defensor Sporting"
$ws.Range("G120").Value = "Danubio"
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 2
$ws.Range("J120").Value = 51
$ws.Range("K120").Value = 1.8
$ws.Range("L120").Value = 3.6
$ws.Range("M120").Value = 4.2
$ws.Range("N120").Value = 1.8
$ws.Range("O120").Value = 3.6
$ws.Range("P120").Value = 4.2
$ws.Range("Q120").Value = -0.75
$ws.Range("R120").Value = 2.05
$ws.Range("S120").Value = 1.8
$ws.Range("T120").Value = 2.25
$ws.Range("U120").Value = 1.85
$ws.Range("V120").Value = 2
$ws.Range("W120").Value = -1
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = 3.2
$ws.Range("Z120").Value = -1
$ws.Range("AA120").Value = 0.8
$ws.Range("AB120").Value = -0.5
$ws.Range("AC120").Value = 0.5

# --- Update rows 169-171 in place (existing future-match rows get refreshed odds) ---
# Row 169
$ws.Range("A169").Value = 167
$ws.Range("B169").Value = 8014133
$ws.Range("C169").Value = "Uruguay Primera División"
$ws.Range("D169").Value = "Uruguay Apertura"
$ws.Range("E169").Value = 45388.41666666666
$ws.Range("F169").Value = "CA River Plate"
$ws.Range("G169").Value = "Montevideo Wanderers"
$ws.Range("K169").Value = 2.5
$ws.Range("L169").Value = 2.9
$ws.Range("M169").Value = 2.9
$ws.Range("N169").Value = 2.5
$ws.Range("O169").Value = 2.9
$ws.Range("P169").Value = 2.875
$ws.Range("Q169").Value = 0
$ws.Range("R169").Value = 1.8
$ws.Range("S169").Value = 2.05
$ws.Range("T169").Value = 2
$ws.Range("U169").Value = 1.85
$ws.Range("V169").Value = 2
$ws.Range("W169").Value = 0
$ws.Range("X169").Value = 0
$ws.Range("Y169").Value = 0
$ws.Range("Z169").Value = 0
$ws.Range("AA169").Value = 0

# Row 170
$ws.Range("A170").Value = 168
$ws.Range("B170").Value = 8014044
$ws.Range("C170").Value = "Uruguay Primera División"
$ws.Range("D170").Value = "Uruguay Apertura"
$ws.Range("E170").Value = 45388.52083333334
$ws.Range("F170").Value = "Racing Club de Montevideo"
$ws.Range("G170").Value = "Miramar Misiones"
$ws.Range("K170").Value = 1.909
$ws.Range("L170").Value = 3.3
$ws.Range("M170").Value = 4
$ws.Range("N170").Value = 1.85
$ws.Range("O170").Value = 3.5
$ws.Range("P170").Value = 4
$ws.Range("Q170").Value = -0.5
$ws.Range("R170").Value = 1.9
$ws.Range("S170").Value = 1.95
$ws.Range("T170").Value = 2.5
$ws.Range("U170").Value = 2.025
$ws.Range("V170").Value = 1.825
$ws.Range("W170").Value = 0
$ws.Range("X170").Value = 0
$ws.Range("Y170").Value = 0
$ws.Range("Z170").Value = 0
$ws.Range("AA170").Value = 0

# Row 171
$ws.Range("A171").Value = 169
$ws.Range("B171").Value = 8014043
$ws.Range("C171").Value = "Uruguay Primera División"
$ws.Range("D171").Value = "Uruguay Apertura"
$ws.Range("E171").Value = 45388.625
$ws.Range("F171").Value = "Danubio"
$ws.Range("G171").Value = "Club Atletico Progreso"
$ws.Range("K171").Value = 2.375
$ws.Range("L171").Value = 3.1
$ws.Range("M171").Value = 3
$ws.Range("N171").Value = 2.4
$ws.Range("O171").Value = 3
$ws.Range("P171").Value = 3
$ws.Range("Q171").Value = -0.25
$ws.Range("R171").Value = 2.1
$ws.Range("S171").Value = 1.775
$ws.Range("T171").Value = 2.25
$ws.Range("U171").Value = 2.025
$ws.Range("V171").Value = 1.825
$ws.Range("W171").Value = 0
$ws.Range("X171").Value = 0
$ws.Range("Y171").Value = 0
$ws.Range("Z171").Value = 0
$ws.Range("AA171").Value = 0

# --- Add 2 new future-match rows (172, 173), copying formats (A,E) from row 171 template ---
$ws.Range("A171").Copy()
$ws.Range("A172").PasteSpecial(-4122)
$ws.Range("E171").Copy()
$ws.Range("E172").PasteSpecial(-4122)
# Row 172
$ws.Range("A172").Value = 170
$ws.Range("B172").Value = 8014089
$ws.Range("C172").Value = "Uruguay Primera División"
$ws.Range("D172").Value = "Uruguay Apertura"
$ws.Range("E172").Value = 45388.75
$ws.Range("F172").Value = "Nacional De Football"
$ws.Range("G172").Value = "Cerro Largo"
$ws.Range("K172").Value = 1.5
$ws.Range("L172").Value = 4
$ws.Range("M172").Value = 6
$ws.Range("N172").Value = 1.55
$ws.Range("O172").Value = 4
$ws.Range("P172").Value = 5.5
$ws.Range("Q172").Value = -1
$ws.Range("R172").Value = 2.05
$ws.Range("S172").Value = 1.8
$ws.Range("T172").Value = 2.25
$ws.Range("U172").Value = 1.9
$ws.Range("V172").Value = 1.95
$ws.Range("W172").Value = 0
$ws.Range("X172").Value = 0
$ws.Range("Y172").Value = 0
$ws.Range("Z172").Value = 0
$ws.Range("AA172").Value = 0

$ws.Range("A171").Copy()
$ws.Range("A173").PasteSpecial(-4122)
$ws.Range("E171").Copy()
$ws.Range("E173").PasteSpecial(-4122)
# Row 173
$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 8014091
$ws.Range("C173").Value = "Uruguay Primera División"
$ws.Range("D173").Value = "Uruguay Apertura"
$ws.Range("E173").Value = 45389.70833333334
$ws.Range("F173").Value = "Deportivo Maldonado"
$ws.Range("G173").Value = "Penarol"
$ws.Range("K173").Value = 5
$ws.Range("L173").Value = 3.75
$ws.Range("M173").Value = 1.615
$ws.Range("N173").Value = 4.2
$ws.Range("O173").Value = 3.75
$ws.Range("P173").Value = 1.7
$ws.Range("Q173").Value = 0.75
$ws.Range("R173").Value = 1.9
$ws.Range("S173").Value = 1.95
$ws.Range("T173").Value = 2.5
$ws.Range("U173").Value = 2.025
$ws.Range("V173").Value = 1.825
$ws.Range("W173").Value = 0
$ws.Range("X173").Value = 0
$ws.Range("Y173").Value = 0
$ws.Range("Z173").Value = 0
$ws.Range("AA173").Value = 0

Write-Host "Edit complete"